# Auto-generated edit script applying the Cactuar_Profits.xlsx data refresh
# (per-cell numeric updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets,
# matching the scheduled-runner commit diff).

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H41").Value = 2009
$ws.Range("I41").Value = 48
$ws.Range("J41").Value = 3098.4443
$ws.Range("K41").Value = 48
$ws.Range("L41").Value = 3098.4443
$ws.Range("M41").Value = 392
$ws.Range("N41").Value = -3978.4443
$ws.Range("H43").Value = 2415158.8
$ws.Range("I43").Value = 4222528
$ws.Range("J43").Value = 5333
$ws.Range("K43").Value = 4222528
$ws.Range("L43").Value = 5333
$ws.Range("M43").Value = -4222459
$ws.Range("N43").Value = -5471
$ws.Range("H64").Value = 26351128
$ws.Range("J64").Value = 55559884
$ws.Range("L64").Value = 55559884
$ws.Range("N64").Value = -55560380
$ws.Range("H67").Value = 26351128
$ws.Range("J67").Value = 55559884
$ws.Range("L67").Value = 55559884
$ws.Range("N67").Value = -55561600
$ws.Range("H70").Value = 4419.077
$ws.Range("I70").Value = 1599.2
$ws.Range("J70").Value = 6181.5
$ws.Range("K70").Value = 4797.6
$ws.Range("L70").Value = 18544.5
$ws.Range("M70").Value = -4527.6
$ws.Range("N70").Value = -19084.5
$ws.Range("H73").Value = 4419.077
$ws.Range("I73").Value = 1599.2
$ws.Range("J73").Value = 6181.5
$ws.Range("K73").Value = 4797.6
$ws.Range("L73").Value = 18544.5
$ws.Range("M73").Value = -3861.6
$ws.Range("N73").Value = -20416.5
$ws.Range("H88").Value = 11773137
$ws.Range("I88").Value = 50002110
$ws.Range("J88").Value = 10375.615
$ws.Range("K88").Value = 50002110
$ws.Range("L88").Value = 10375.615
$ws.Range("M88").Value = -50001704
$ws.Range("N88").Value = -11187.615
$ws.Range("H91").Value = 11773137
$ws.Range("I91").Value = 50002110
$ws.Range("J91").Value = 10375.615
$ws.Range("K91").Value = 50002110
$ws.Range("L91").Value = 10375.615
$ws.Range("M91").Value = -50000706
$ws.Range("N91").Value = -13183.615
$ws.Range("H98").Value = 1020.375
$ws.Range("I98").Value = 1020.375
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1020.375
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 477.625
$ws.Range("N98").ClearContents()
$ws.Range("H107").Value = 253.9375
$ws.Range("I107").Value = 256
$ws.Range("K107").Value = 256
$ws.Range("M107").Value = 1664
$ws.Range("H113").Value = 3633.682
$ws.Range("I113").Value = 2794.5
$ws.Range("J113").Value = 4333
$ws.Range("K113").Value = 2794.5
$ws.Range("L113").Value = 4333
$ws.Range("M113").Value = 459.5
$ws.Range("N113").Value = -10841
$ws.Range("H116").Value = 34728224
$ws.Range("I116").Value = 20067656
$ws.Range("J116").Value = 166673330
$ws.Range("K116").Value = 20067656
$ws.Range("L116").Value = 166673330
$ws.Range("M116").Value = -20064214
$ws.Range("N116").Value = -166680214
$ws.Range("H122").Value = 1020.375
$ws.Range("I122").Value = 1020.375
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3061.125
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -611.125
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 15813.635
$ws.Range("I132").Value = 7092.64
$ws.Range("J132").Value = 23888.629
$ws.Range("K132").Value = 21277.92
$ws.Range("L132").Value = 71665.887
$ws.Range("M132").Value = -18747.92
$ws.Range("N132").Value = -76725.887
$ws.Range("H137").Value = 11079741
$ws.Range("I137").Value = 716548
$ws.Range("J137").Value = 19614136
$ws.Range("K137").Value = 2149644
$ws.Range("L137").Value = 58842408
$ws.Range("M137").Value = -2147094
$ws.Range("N137").Value = -58847508
$ws.Range("H138").Value = 2161.1099
$ws.Range("J138").Value = 2619.3684
$ws.Range("L138").Value = 7858.1052
$ws.Range("N138").Value = -18138.1052

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 584347.25
$ws.Range("I2").Value = 761635.8
$ws.Range("J2").Value = 1827.5714
$ws.Range("K2").Value = 761635.8
$ws.Range("L2").Value = 1827.5714
$ws.Range("M2").Value = -761522.8
$ws.Range("N2").Value = -2053.5714
$ws.Range("H63").Value = 4098.5
$ws.Range("I63").Value = 1950
$ws.Range("K63").Value = 1950
$ws.Range("M63").Value = -1264
$ws.Range("H66").Value = 4098.5
$ws.Range("I66").Value = 1950
$ws.Range("K66").Value = 9750
$ws.Range("M66").Value = -6318
$ws.Range("H116").Value = 584347.25
$ws.Range("I116").Value = 761635.8
$ws.Range("J116").Value = 1827.5714
$ws.Range("K116").Value = 761635.8
$ws.Range("L116").Value = 1827.5714
$ws.Range("M116").Value = -759341.8
$ws.Range("N116").Value = -6415.5714
$ws.Range("H132").Value = 12061.692
$ws.Range("I132").Value = 14678.452
$ws.Range("J132").Value = 7283.2607
$ws.Range("K132").Value = 44035.356
$ws.Range("L132").Value = 21849.7821
$ws.Range("M132").Value = -41505.356
$ws.Range("N132").Value = -26909.7821

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 584347.25
$ws.Range("I3").Value = 761635.8
$ws.Range("J3").Value = 1827.5714
$ws.Range("K3").Value = 761635.8
$ws.Range("L3").Value = 1827.5714
$ws.Range("M3").Value = -761521.8
$ws.Range("N3").Value = -2055.5714
$ws.Range("H20").Value = 9175.214
$ws.Range("J20").Value = 9602
$ws.Range("L20").Value = 9602
$ws.Range("N20").Value = -10096

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 2154.0908
$ws.Range("I16").Value = 2021.6666
$ws.Range("K16").Value = 2021.6666
$ws.Range("M16").Value = -1734.6666
$ws.Range("H31").Value = 5895.8086
$ws.Range("I31").Value = 2289.0625
$ws.Range("K31").Value = 2289.0625
$ws.Range("M31").Value = -1994.0625
$ws.Range("H34").Value = 5895.8086
$ws.Range("I34").Value = 2289.0625
$ws.Range("K34").Value = 2289.0625
$ws.Range("M34").Value = -2087.0625
$ws.Range("H39").Value = 5250
$ws.Range("I39").Value = 5250
$ws.Range("K39").Value = 5250
$ws.Range("M39").Value = -4859
$ws.Range("H49").Value = 5250
$ws.Range("I49").Value = 5250
$ws.Range("K49").Value = 5250
$ws.Range("M49").Value = -5068
$ws.Range("H99").Value = 19455.637
$ws.Range("I99").Value = 110012
$ws.Range("J99").Value = 10400
$ws.Range("K99").Value = 110012
$ws.Range("L99").Value = 10400
$ws.Range("M99").Value = -108514
$ws.Range("N99").Value = -13396
$ws.Range("H105").Value = 4546764.5
$ws.Range("J105").Value = 1833.3334
$ws.Range("L105").Value = 1833.3334
$ws.Range("N105").Value = -5327.3334
$ws.Range("H113").Value = 2154.0908
$ws.Range("I113").Value = 2021.6666
$ws.Range("K113").Value = 2021.6666
$ws.Range("M113").Value = 148.3334
$ws.Range("H126").Value = 19455.637
$ws.Range("I126").Value = 110012
$ws.Range("J126").Value = 10400
$ws.Range("K126").Value = 330036
$ws.Range("L126").Value = 31200
$ws.Range("M126").Value = -327566
$ws.Range("N126").Value = -36140
$ws.Range("H132").Value = 15163263
$ws.Range("I132").Value = 20848832
$ws.Range("J132").Value = 1749
$ws.Range("K132").Value = 62546496
$ws.Range("L132").Value = 5247
$ws.Range("M132").Value = -62543966
$ws.Range("N132").Value = -10307
$ws.Range("H141").Value = 92872.7
$ws.Range("J141").Value = 106264.93
$ws.Range("L141").Value = 106264.93
$ws.Range("N141").Value = -116624.93

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H11").Value = 136668800
$ws.Range("I11").Value = 945
$ws.Range("J11").Value = 205002720
$ws.Range("K11").Value = 2835
$ws.Range("L11").Value = 615008160
$ws.Range("M11").Value = -2695
$ws.Range("N11").Value = -615008440
$ws.Range("H48").Value = 9999.200000000001
$ws.Range("J48").Value = 9999.200000000001
$ws.Range("L48").Value = 29997.6
$ws.Range("N48").Value = -30497.6
$ws.Range("H131").Value = 12858582
$ws.Range("J131").Value = 13260384
$ws.Range("L131").Value = 39781152
$ws.Range("N131").Value = -39791232

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 1693024.8
$ws.Range("I80").Value = 5542583
$ws.Range("K80").Value = 5542583
$ws.Range("M80").Value = -5541585
$ws.Range("H83").Value = 1693024.8
$ws.Range("I83").Value = 5542583
$ws.Range("K83").Value = 27712915
$ws.Range("M83").Value = -27707923

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 884
$ws.Range("I22").Value = 753.36365
$ws.Range("K22").Value = 753.36365
$ws.Range("M22").Value = -458.36365
$ws.Range("H27").Value = 884
$ws.Range("I27").Value = 753.36365
$ws.Range("K27").Value = 753.36365
$ws.Range("M27").Value = -646.36365
$ws.Range("H46").Value = 6977
$ws.Range("I46").Value = 4600
$ws.Range("J46").Value = 7116.8237
$ws.Range("K46").Value = 4600
$ws.Range("L46").Value = 7116.8237
$ws.Range("M46").Value = -4412
$ws.Range("N46").Value = -7492.8237
$ws.Range("H68").Value = 3249146.5
$ws.Range("I68").Value = 3789336.2
$ws.Range("J68").Value = 8008
$ws.Range("K68").Value = 3789336.2
$ws.Range("L68").Value = 8008
$ws.Range("M68").Value = -3788587.2
$ws.Range("N68").Value = -9506
$ws.Range("H71").Value = 3249146.5
$ws.Range("I71").Value = 3789336.2
$ws.Range("J71").Value = 8008
$ws.Range("K71").Value = 18946681
$ws.Range("L71").Value = 40040
$ws.Range("M71").Value = -18942937
$ws.Range("N71").Value = -47528
$ws.Range("H93").Value = 2080.3447
$ws.Range("I93").Value = 2225.32
$ws.Range("K93").Value = 2225.32
$ws.Range("M93").Value = -977.3200000000002
$ws.Range("H122").Value = 90914210
$ws.Range("I122").Value = 142861460
$ws.Range("J122").Value = 6499.5
$ws.Range("K122").Value = 428584380
$ws.Range("L122").Value = 19498.5
$ws.Range("M122").Value = -428581930
$ws.Range("N122").Value = -24398.5
$ws.Range("H132").Value = 4624.276
$ws.Range("I132").Value = 4601.125
$ws.Range("J132").Value = 4652.769
$ws.Range("K132").Value = 13803.375
$ws.Range("L132").Value = 13958.307
$ws.Range("M132").Value = -11273.375
$ws.Range("N132").Value = -19018.307
$ws.Range("H136").Value = 5085.7144
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
